# New crime data collected - update weekly CompStat report figures
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$intFmt = "#,##0"
$pctFmt = "#,##0.0;""-""#,##0.0"

# --- Header text updates (report volume/number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 30   Number  19"
$ws.Range("C9").Value = "Report Covering the Week  5/8/2023  Through  5/14/2023"

# --- Row 15 (Rape) ---
$ws.Range("D15").Value = 1
$ws.Range("D15").NumberFormat = $intFmt
$ws.Range("E15").Value = -100
$ws.Range("E15").NumberFormat = $pctFmt
$ws.Range("G15").Value = 1
$ws.Range("G15").NumberFormat = $intFmt
$ws.Range("H15").Value = 0
$ws.Range("H15").NumberFormat = $pctFmt
$ws.Range("J15").Value = 6
$ws.Range("K15").Value = -16.666666666666

# --- Row 16 (Robbery) ---
$ws.Range("C16").Value = 2
$ws.Range("C16").NumberFormat = $intFmt
$ws.Range("D16").Value = 4
$ws.Range("D16").NumberFormat = $intFmt
$ws.Range("E16").Value = -50
$ws.Range("E16").NumberFormat = $pctFmt
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 8
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 64
$ws.Range("J16").Value = 61
$ws.Range("K16").Value = 4.918032786885
$ws.Range("L16").Value = 10.344827586206
$ws.Range("M16").Value = -24.705882352941
$ws.Range("N16").Value = -61.676646706586

# --- Row 17 (Fel. Assault) ---
$ws.Range("C17").Value = 10
$ws.Range("E17").Value = 66.666666666666
$ws.Range("F17").Value = 22
$ws.Range("G17").Value = 16
$ws.Range("H17").Value = 37.5
$ws.Range("I17").Value = 110
$ws.Range("J17").Value = 85
$ws.Range("K17").Value = 29.411764705882
$ws.Range("L17").Value = 50.684931506849
$ws.Range("M17").Value = 66.666666666666
$ws.Range("N17").Value = 37.5

# --- Row 18 (Burglary) ---
$ws.Range("D18").Value = 3
$ws.Range("E18").Value = 0
$ws.Range("G18").Value = 12
$ws.Range("H18").Value = 25
$ws.Range("I18").Value = 56
$ws.Range("J18").Value = 40
$ws.Range("K18").Value = 40
$ws.Range("L18").Value = 47.368421052631
$ws.Range("M18").Value = -41.052631578947
$ws.Range("N18").Value = -82.389937106918

# --- Row 19 (Gr. Larceny) ---
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 150
$ws.Range("F19").Value = 59
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 78.787878787878
$ws.Range("I19").Value = 242
$ws.Range("J19").Value = 197
$ws.Range("K19").Value = 22.842639593908
$ws.Range("L19").Value = 68.055555555555
$ws.Range("M19").Value = 63.513513513513
$ws.Range("N19").Value = 55.128205128205

# --- Row 20 (G.L.A.) ---
$ws.Range("C20").Value = 11
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 450
$ws.Range("F20").Value = 51
$ws.Range("G20").Value = 18
$ws.Range("H20").Value = 183.333333333333
$ws.Range("I20").Value = 196
$ws.Range("J20").Value = 92
$ws.Range("K20").Value = 113.04347826087
$ws.Range("L20").Value = 192.537313432836
$ws.Range("M20").Value = 221.311475409836
$ws.Range("N20").Value = -72.549019607843

# --- Row 21 (TOTAL) ---
$ws.Range("C21").Value = 41
$ws.Range("E21").Value = 86.363636363636
$ws.Range("F21").Value = 156
$ws.Range("G21").Value = 88
$ws.Range("H21").Value = 77.272727272727
$ws.Range("I21").Value = 674
$ws.Range("J21").Value = 481
$ws.Range("K21").Value = 40.124740124740
$ws.Range("L21").Value = 75.064935064935
$ws.Range("M21").Value = 45.887445887445
$ws.Range("N21").Value = -53.388658367911

# --- Row 22 (Transit) ---
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50
$ws.Range("J22").Value = 5
$ws.Range("K22").Value = 40

# --- Row 23 (Housing) ---
$ws.Range("D23").Value = 1
$ws.Range("D23").NumberFormat = $intFmt
$ws.Range("E23").Value = 100
$ws.Range("E23").NumberFormat = $pctFmt
$ws.Range("F23").Value = 5
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 66.666666666666
$ws.Range("I23").Value = 20
$ws.Range("J23").Value = 13
$ws.Range("K23").Value = 53.846153846153
$ws.Range("L23").Value = 66.666666666666
$ws.Range("M23").Value = 53.846153846153

# --- Row 24 (Petit Larceny) ---
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = 11.111111111111
$ws.Range("F24").Value = 89
$ws.Range("G24").Value = 106
$ws.Range("H24").Value = -16.037735849056
$ws.Range("I24").Value = 497
$ws.Range("J24").Value = 414
$ws.Range("K24").Value = 20.048309178744
$ws.Range("L24").Value = 53.395061728395
$ws.Range("M24").Value = 5.744680851063

# --- Row 25 (Misd. Assault) ---
$ws.Range("C25").Value = 11
$ws.Range("D25").Value = 8
$ws.Range("E25").Value = 37.5
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 36
$ws.Range("H25").Value = 22.222222222222
$ws.Range("I25").Value = 188
$ws.Range("J25").Value = 169
$ws.Range("K25").Value = 11.242603550295
$ws.Range("L25").Value = 27.891156462585
$ws.Range("M25").Value = 33.333333333333

# --- Row 26 (UCR Rape*) ---
$ws.Range("D26").Value = 3
$ws.Range("D26").NumberFormat = $intFmt
$ws.Range("E26").Value = -100
$ws.Range("E26").NumberFormat = $pctFmt
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = -50
$ws.Range("J26").Value = 15
$ws.Range("K26").Value = -33.333333333333

# --- Row 27 (Other Sex Crimes) ---
$ws.Range("F27").Value = 5
$ws.Range("G27").Value = 6
$ws.Range("H27").Value = -16.666666666666
$ws.Range("J27").Value = 21
$ws.Range("K27").Value = -28.571428571428
$ws.Range("L27").Value = -21.052631578947
